# ---------------------------------------------------------------------------
# UTAH_2024.xlsx cleanup
#  1. Rename header row (row 1) from Spanish display labels to snake_case
#     column keys used by the data pipeline.
#  2. Normalize Spanish grammatical connector words ("de", "del", "la",
#     "las", "los", "el", "y") inside state/municipality names to Title
#     Case (e.g. "Pabellón de Arteaga" -> "Pabellón De Arteaga").
#  3. Fix the grand-total label in A1091 from the all-caps "TOTAL" to
#     "Total" (matching the style used by every subtotal row above it).
#  4. Drop the trailing free-text metadata/footnote rows (1093-1097) that
#     aren't part of the tabular dataset, and let the sheet's used range /
#     dimension shrink back down to A1:D1091.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row: Spanish labels -> snake_case machine-readable keys ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case the lowercase Spanish connector words in names ---
$renames = @{
    "B7" = "Pabellón De Arteaga"
    "B8" = "Rincón De Romos"
    "B29" = "Amatenango De La Frontera"
    "B33" = "Benemérito De Las Américas"
    "B40" = "Chiapa De Corzo"
    "B61" = "Ocozocoautla De Espinosa"
    "B66" = "San Cristóbal De Las Casas"
    "B98" = "Guadalupe Y Calvo"
    "B102" = "Hidalgo Del Parral"
    "B116" = "San Francisco Del Oro"
    "B121" = "Valle De Zaragoza"
    "A123" = "Ciudad De México"
    "B126" = "Cuajimalpa De Morelos"
    "A140" = "Coahuila De Zaragoza"
    "B158" = "Villa De Álvarez"
    "B173" = "Nombre De Dios"
    "B178" = "Pánuco De Coronado"
    "A186" = "Estado De México"
    "B186" = "Acambay De Ruíz Castañeda"
    "B188" = "Almoloya De Alquisiras"
    "B189" = "Almoloya De Juárez"
    "B193" = "Atizapán De Zaragoza"
    "B198" = "Chapa De Mota"
    "B201" = "Coacalco De Berriozábal"
    "B207" = "Ecatepec De Morelos"
    "B210" = "Ixtapan De La Sal"
    "B220" = "Naucalpan De Juárez"
    "B227" = "San Felipe Del Progreso"
    "B228" = "San José Del Rincón"
    "B229" = "San Martín De Las Pirámides"
    "B239" = "Tenango Del Aire"
    "B240" = "Tenango Del Valle"
    "B247" = "Tlalnepantla De Baz"
    "B251" = "Valle De Bravo"
    "B262" = "Apaseo El Alto"
    "B263" = "Apaseo El Grande"
    "B269" = "Dolores Hidalgo Cuna De La Independencia Nacional"
    "B272" = "Jaral Del Progreso"
    "B278" = "Purísima Del Rincón"
    "B283" = "San Diego De La Unión"
    "B285" = "San Francisco Del Rincón"
    "B287" = "San Luis De La Paz"
    "B288" = "San Miguel De Allende"
    "B289" = "Santa Cruz De Juventino Rosas"
    "B290" = "Silao De La Victoria"
    "B294" = "Valle De Santiago"
    "B299" = "Acapulco De Juárez"
    "B302" = "Ajuchitlán Del Progreso"
    "B303" = "Alcozauca De Guerrero"
    "B307" = "Atenango Del Río"
    "B308" = "Atoyac De Álvarez"
    "B309" = "Ayutla De Los Libres"
    "B312" = "Buenavista De Cuéllar"
    "B313" = "Chilapa De Álvarez"
    "B314" = "Chilpancingo De Los Bravo"
    "B315" = "Coahuayutla De José María Izazaga"
    "B320" = "Coyuca De Benítez"
    "B321" = "Coyuca De Catalán"
    "B324" = "Cuetzala Del Progreso"
    "B331" = "Huitzuco De Los Figueroa"
    "B332" = "Iguala De La Independencia"
    "B335" = "La Unión De Isidoro Montes De Oca"
    "B341" = "Mártir De Cuilapan"
    "B352" = "Taxco De Alarcón"
    "B355" = "Tepecoacuilco De Trujano"
    "B357" = "Tixtla De Guerrero"
    "B358" = "Tlalixtaquilla De Maldonado"
    "B359" = "Tlapa De Comonfort"
    "B361" = "Técpan De Galeana"
    "B364" = "Zihuatanejo De Azueta"
    "B370" = "Agua Blanca De Iturbide"
    "B374" = "Atotonilco De Tula"
    "B375" = "Atotonilco El Grande"
    "B378" = "Cuautepec De Hinojosa"
    "B382" = "Huasca De Ocampo"
    "B385" = "Huejutla De Reyes"
    "B389" = "Mineral Del Monte"
    "B390" = "Mixquiahuala De Juárez"
    "B391" = "Molango De Escamilla"
    "B392" = "Omitlán De Juárez"
    "B393" = "Pachuca De Soto"
    "B396" = "Santiago Tulantepec De Lugo Guerrero"
    "B399" = "Tenango De Doria"
    "B401" = "Tepehuacán De Guerrero"
    "B402" = "Tepeji Del Río De Ocampo"
    "B403" = "Tezontepec De Aldama"
    "B408" = "Tula De Allende"
    "B409" = "Tulancingo De Bravo"
    "B412" = "Zacualtipán De Ángeles"
    "B415" = "Acatlán De Juárez"
    "B418" = "Atotonilco El Alto"
    "B420" = "Autlán De Navarro"
    "B430" = "Cuautitlán De García Barragán"
    "B434" = "Encarnación De Díaz"
    "B439" = "Huejuquilla El Alto"
    "B441" = "Ixtlahuacán De Los Membrillos"
    "B445" = "Jilotlán De Los Dolores"
    "B450" = "La Manzanilla De La Paz"
    "B451" = "Lagos De Moreno"
    "B458" = "San Diego De Alejandría"
    "B460" = "San Juanito De Escobedo"
    "B462" = "San Miguel El Alto"
    "B464" = "Santa María De Los Ángeles"
    "B467" = "Talpa De Allende"
    "B468" = "Tamazula De Gordiano"
    "B473" = "Teocuitatlán De Corona"
    "B474" = "Tepatitlán De Morelos"
    "B477" = "Tizapán El Alto"
    "B478" = "Tlajomulco De Zúñiga"
    "B486" = "Unión De San Antonio"
    "B487" = "Unión De Tula"
    "B488" = "Valle De Juárez"
    "B493" = "Yahualica De González Gallo"
    "B494" = "Zacoalco De Torres"
    "B497" = "Zapotitlán De Vadillo"
    "B499" = "Zapotlán Del Rey"
    "B500" = "Zapotlán El Grande"
    "A502" = "Michoacán De Ocampo"
    "B517" = "Coalcomán De Vázquez Pallares"
    "B595" = "Coatlán Del Río"
    "B605" = "Puente De Ixtla"
    "B610" = "Tetela Del Volcán"
    "B611" = "Tlaltizapán De Zapata"
    "B617" = "Zacualpan De Amilpas"
    "B621" = "Amatlán De Cañas"
    "B622" = "Bahía De Banderas"
    "B627" = "Ixtlán Del Río"
    "B632" = "Santa María Del Oro"
    "B641" = "San Nicolás De Los Garza"
    "B644" = "Acatlán De Pérez Figueroa"
    "B649" = "Cuyamecalco Villa De Zaragoza"
    "B652" = "Guevea De Humboldt"
    "B653" = "Heroica Ciudad De Huajuapan De León"
    "B654" = "Heroica Ciudad De Juchitán De Zaragoza"
    "B655" = "Huajuapan De León"
    "B656" = "Huautla De Jiménez"
    "B657" = "Ixtlán De Juárez"
    "B663" = "Mártires De Tacubaya"
    "B664" = "Oaxaca De Juárez"
    "B665" = "Ocotlán De Morelos"
    "B666" = "Putla Villa De Guerrero"
    "B671" = "San Antonio De La Cal"
    "B674" = "San Dionisio Del Mar"
    "B679" = "San Francisco Del Mar"
    "B683" = "San José Del Progreso"
    "B686" = "San Juan Bautista Lo De Soto"
    "B706" = "San Pedro Y San Pablo Teposcolula"
    "B717" = "Santa María Jalapa Del Marqués"
    "B742" = "Tataltepec De Valdés"
    "B743" = "Teotitlán De Flores Magón"
    "B744" = "Tlacolula De Matamoros"
    "B745" = "Villa De Etla"
    "B746" = "Villa De Tamazulápam Del Progreso"
    "B747" = "Villa De Tututepec"
    "B748" = "Villa De Tututepec De Melchor Ocampo"
    "B749" = "Villa De Zaachila"
    "B750" = "Zimatlán De Álvarez"
    "B765" = "Chalchicomula De Sesma"
    "B778" = "Cuayuca De Andrade"
    "B785" = "Huehuetlán El Chico"
    "B786" = "Huehuetlán El Grande"
    "B789" = "Izúcar De Matamoros"
    "B794" = "Los Reyes De Juárez"
    "B803" = "Palmar De Bravo"
    "B813" = "San Salvador El Seco"
    "B816" = "Tecali De Herrera"
    "B820" = "Tepanco De López"
    "B821" = "Tepatlaxco De Hidalgo"
    "B823" = "Tepexi De Rodríguez"
    "B824" = "Tetela De Ocampo"
    "B828" = "Tlacotepec De Benito Juárez"
    "B844" = "Amealco De Bonfil"
    "B845" = "Cadereyta De Montes"
    "B846" = "Jalpan De Serra"
    "B850" = "San Juan Del Río"
    "B864" = "Mexquitic De Carmona"
    "B874" = "Villa De Arriaga"
    "B875" = "Villa De Ramos"
    "B876" = "Villa De Reyes"
    "B913" = "Nacozari De García"
    "B920" = "San Pedro De La Cueva"
    "B930" = "Jalpa De Méndez"
    "B946" = "Contla De Juan Cuamatzi"
    "B950" = "Ixtacuixtla De Mariano Matamoros"
    "B953" = "San Pablo Del Monte"
    "B954" = "Sanctórum De Lázaro Cárdenas"
    "A962" = "Veracruz De Ignacio De La Llave"
    "B967" = "Amatlán De Los Reyes"
    "B971" = "Boca Del Río"
    "B984" = "Cosamaloapan De Carpio"
    "B985" = "Cosautlán De Carvajal"
    "B994" = "Hueyapan De Ocampo"
    "B995" = "Ignacio De La Llave"
    "B996" = "Ixhuatlán De Madero"
    "B1005" = "Martínez De La Torre"
    "B1010" = "Nanchital De Lázaro Cárdenas Del Río"
    "B1016" = "Paso Del Macho"
    "B1018" = "Poza Rica De Hidalgo"
    "B1026" = "Sayula De Alemán"
    "B1028" = "Soledad De Doblado"
    "B1063" = "Concepción Del Oro"
    "B1074" = "Noria De Ángeles"
    "B1083" = "Teúl De González Ortega"
    "B1084" = "Tlaltenango De Sánchez Román"
    "B1087" = "Villa De Cos"
}

foreach ($ref in $renames.Keys) {
    $ws.Range($ref).Value = $renames[$ref]
}

# --- 3. Fix the grand-total row label ---
$ws.Range("A1091").Value = "Total"

# --- 4. Remove the trailing metadata/footnote rows (1093-1097) ---
$ws.Rows.Item(1093).Resize(5).Delete()
